# Add the 11 May 2020 (11 Mayis 2020) COVID-19 Turkey daily data row.
# Source data was sheet "data" backed by table "Table3" (A1:E60); a new
# row (61) is appended with date serial 43962 (2020-05-11) and that
# day's test/case/death/recovered counts, then the table is resized to
# include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (60) down into the
# new row (61) so the new cells pick up the same number formats /
# alignment (date format for column A, left-aligned general for the
# rest) as the rest of the table.
$ws.Range("A60:E60").Copy() | Out-Null
$ws.Range("A61:E61").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New day's values: date, test, case, death, recovered.
$ws.Range("A61").Value2 = 43962
$ws.Range("B61").Value2 = 32722
$ws.Range("C61").Value2 = 1114
$ws.Range("D61").Value2 = 55
$ws.Range("E61").Value2 = 3089

# Grow the table (ListObject) so it covers the newly added row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E61"))

# Match the updated active-cell selection recorded in the saved file.
$ws.Range("E60").Select() | Out-Null
